$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Column D now carries a combined "date - description" label for the
# single homework/attendance session instead of a second raw date.
$ws.Range("D1").Value = "2025-02-23 - HW1"

# --- Row 2 (student record was swapped for a new student) -------------
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "650610759"
$ws.Range("B2").Value = "earn"
$ws.Range("C2").Value = "earn@example.com"
# D2 ("มาเรียน") is unchanged

# --- Row 3 (attendance result updated) ---------------------------------
$ws.Range("D3").Value = "ขาดเรียน"

# --- Drop the now-unused second attendance/date column (E) -------------
# This shifts nothing else; used range becomes A1:D3.
$ws.Columns("E:E").Delete()

# --- Widen the email / date-description columns to fit the new text ----
$ws.Range("C1").EntireColumn.ColumnWidth = 18.83203125
$ws.Range("D1").EntireColumn.ColumnWidth = 18.83203125
